$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.486.71"
$ws.Range("E2").Value = "'  -2.19%  "

# Row 3
$ws.Range("D3").Value = "'1.748.35"
$ws.Range("E3").Value = "'  -2.40%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.11%  "

# Row 5
$ws.Range("D5").Value = "'324.59"
$ws.Range("E5").Value = "'  +0.12%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.02%  "

# Row 7
$ws.Range("D7").Value = "'0.4471"
$ws.Range("E7").Value = "'  +4.01%  "

# Row 8
$ws.Range("D8").Value = "'0.3596"
$ws.Range("E8").Value = "'  -1.03%  "

# Row 9
$ws.Range("D9").Value = "'0.07494"
$ws.Range("E9").Value = "'  -0.50%  "

# Row 10
$ws.Range("D10").Value = "'41.92"
$ws.Range("E10").Value = "'  -5.87%  "

# Row 11
$ws.Range("D11").Value = "'1.093"
$ws.Range("E11").Value = "'  -1.83%  "

# Row 12
$ws.Range("E12").Value = "'  +0.00%  "

# Row 13
$ws.Range("D13").Value = "'20.67"
$ws.Range("E13").Value = "'  -5.00%  "

# Row 14
$ws.Range("D14").Value = "'6.021"
$ws.Range("E14").Value = "'  -2.43%  "

# Row 15
$ws.Range("D15").Value = "'7.133"
$ws.Range("E15").Value = "'  -3.03%  "

# Row 16
$ws.Range("D16").Value = "'1.748.80"
$ws.Range("E16").Value = "'  -1.71%  "

# Row 17
$ws.Range("D17").Value = "'93.52"
$ws.Range("E17").Value = "'  +1.82%  "

# Row 18
$ws.Range("D18").Value = "'0.00001061"
$ws.Range("E18").Value = "'  -0.59%  "

# Row 19
$ws.Range("D19").Value = "'0.06381"
$ws.Range("E19").Value = "'  +0.57%  "

# Row 20
$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = "'  +0.03%  "

# Row 21
$ws.Range("D21").Value = "'16.78"
$ws.Range("E21").Value = "'  -2.89%  "

# Row 22
$ws.Range("D22").Value = "'5.846"
$ws.Range("E22").Value = "'  -1.99%  "

# Row 23
$ws.Range("D23").Value = "'27.554.24"

# Row 24
$ws.Range("D24").Value = "'11.24"
$ws.Range("E24").Value = "'  -1.73%  "

# Row 25
$ws.Range("D25").Value = "'2.083"
$ws.Range("E25").Value = "'  -3.54%  "

# Row 26
$ws.Range("D26").Value = "'162.37"
$ws.Range("E26").Value = "'  +1.27%  "

# Row 27
$ws.Range("D27").Value = "'20.45"
$ws.Range("E27").Value = "'  +0.30%  "

# Row 28
$ws.Range("D28").Value = "'1.950.65"
$ws.Range("E28").Value = "'  -1.82%  "

# Row 29
$ws.Range("D29").Value = "'2.083"
$ws.Range("E29").Value = "'  -4.63%  "

# Row 30
$ws.Range("D30").Value = "'125.52"
$ws.Range("E30").Value = "'  -1.11%  "

# Row 31
$ws.Range("D31").Value = "'1.081"
$ws.Range("E31").Value = "'  -7.52%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09078"
$ws.Range("E32").Value = "'  +0.68%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.665"
$ws.Range("E33").Value = "'  +4.63%  "

# Row 34
$ws.Range("D34").Value = "'5.537"
$ws.Range("E34").Value = "'  -3.41%  "

# Row 35
$ws.Range("D35").Value = "'11.95"
$ws.Range("E35").Value = "'  -5.91%  "

# Row 36
$ws.Range("E36").Value = "'  -1.42%  "

# Row 37
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "'0.6376"
$ws.Range("E37").Value = "'  -1.53%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06015"
$ws.Range("E38").Value = "'  -0.93%  "

# Row 39
$ws.Range("E39").Value = "'  -1.61%  "

# Row 40
$ws.Range("D40").Value = "'4.941"
$ws.Range("E40").Value = "'  -3.18%  "

# Row 41
$ws.Range("D41").Value = "'1.205"
$ws.Range("E41").Value = "'  +1.74%  "

# Row 42
$ws.Range("D42").Value = "'1.382"
$ws.Range("E42").Value = "'  -2.22%  "

# Row 43
$ws.Range("E43").Value = "'  -1.53%  "

# Row 44
$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = "'  -3.68%  "

# Row 46
$ws.Range("D46").Value = "'0.5897"
$ws.Range("E46").Value = "'  -1.75%  "

# Row 47
$ws.Range("D47").Value = "'122.67"
$ws.Range("E47").Value = "'  -1.48%  "

# Row 48
$ws.Range("D48").Value = "'1.956"
$ws.Range("E48").Value = "'  -1.88%  "

# Row 49
$ws.Range("D49").Value = "'1.148"
$ws.Range("E49").Value = "'  -0.57%  "

# Row 50
$ws.Range("D50").Value = "'0.06853"
$ws.Range("E50").Value = "'  -1.52%  "

# Row 51
$ws.Range("E51").Value = "'  -3.52%  "
